$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 21-40 contain the "Unresolved Conditions" / "Additional Benefits"
# list strings.  The shared-strings table backing these cells was
# reshuffled, which (since the cells are referenced by index/position)
# results in the values visible in these rows changing to a new order
# while keeping the same overall set of strings.  Reproduce the new
# order by writing out the values directly.

$newValues = @(
    "peeling",
    "peeling and dehydrated (due to use of products to control oil)",
    "rough on some portions of the skin",
    "clogged pores/cell accumulation",
    "pre-disposed to acne",
    "excessive oil in some areas (possibly int he t-zone)",
    "scaling",
    "dry/rough skin in feeling and/or appearance/rough skin in feeling and/or appearance and/or flaky on some portions of skin",
    "inflamed and/or irritated",
    "dry/rough skin in feeling and/or appearance",
    "oily (overproduction of sebum)",
    "dry",
    "tewl",
    "flaking",
    "eczema",
    "enlarged pores",
    "itchy",
    "sensitive",
    "itching and/or redness",
    "chapped"
)

$startRow = 21
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
